$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "66.128.75"
Set-TextValue $ws.Range("E2") "  -5.97%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.174.71"
Set-TextValue $ws.Range("E3") "  -9.55%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.997"
Set-TextValue $ws.Range("E4") "  -0.21%  "

# Row 5
Set-TextValue $ws.Range("D5") "566.82"

# Row 6
Set-TextValue $ws.Range("D6") "147.14"
Set-TextValue $ws.Range("E6") "  -14.87%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  -0.10%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.165.53"
Set-TextValue $ws.Range("E8") "  -9.76%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.535"
Set-TextValue $ws.Range("E9") "  -11.89%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.168"
Set-TextValue $ws.Range("E10") "  -14.18%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.34"
Set-TextValue $ws.Range("E11") "  -11.38%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.491"
Set-TextValue $ws.Range("E12") "  -16.62%  "

# Row 13
Set-TextValue $ws.Range("D13") "38.00"
Set-TextValue $ws.Range("E13") "  -18.00%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.0000238"
Set-TextValue $ws.Range("E14") "  -13.68%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.681.49"
Set-TextValue $ws.Range("E15") "  -9.62%  "

# Row 16
Set-TextValue $ws.Range("D16") "65.966.32"
Set-TextValue $ws.Range("E16") "  -6.12%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.166.93"
Set-TextValue $ws.Range("E17") "  -10.07%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -6.88%  "

# Row 19
Set-TextValue $ws.Range("D19") "527.42"
Set-TextValue $ws.Range("E19") "  -13.88%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.04"
Set-TextValue $ws.Range("E20") "  -16.82%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.85"
Set-TextValue $ws.Range("E21") "  -16.33%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.748"
Set-TextValue $ws.Range("E22") "  -15.13%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.69"
Set-TextValue $ws.Range("E23") "  -14.35%  "

# Row 24
Set-TextValue $ws.Range("D24") "84.23"
Set-TextValue $ws.Range("E24") "  -14.67%  "

# Row 25
Set-TextValue $ws.Range("D25") "13.32"
Set-TextValue $ws.Range("E25") "  -14.81%  "

# Row 26
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  -0.05%  "

# Row 27
Set-TextValue $ws.Range("E27") "  -18.14%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.15"
Set-TextValue $ws.Range("E28") "  -17.07%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.89"
Set-TextValue $ws.Range("E29") "  -13.84%  "

# Row 30
Set-TextValue $ws.Range("D30") "28.70"
Set-TextValue $ws.Range("E30") "  -14.70%  "

# Row 31
Set-TextValue $ws.Range("D31") "2.54"
Set-TextValue $ws.Range("E31") "  -14.94%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.11"
Set-TextValue $ws.Range("E32") "  -14.75%  "

# Row 33
Set-TextValue $ws.Range("D33") "6.40"
Set-TextValue $ws.Range("E33") "  -21.00%  "

# Row 34
Set-TextValue $ws.Range("B34") "NEARProtocol"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "5.60"
Set-TextValue $ws.Range("E34") "  -17.83%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.999"
Set-TextValue $ws.Range("E35") "  -0.11%  "

# Row 36
Set-TextValue $ws.Range("B36") "Bittensor"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D36") "514.17"
Set-TextValue $ws.Range("E36") "  -16.06%  "

# Row 37
Set-TextValue $ws.Range("D37") "52.81"
Set-TextValue $ws.Range("E37") "  -7.21%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.0844"
Set-TextValue $ws.Range("E38") "  -15.86%  "

# Row 39
Set-TextValue $ws.Range("D39") "9.02"
Set-TextValue $ws.Range("E39") "  -16.53%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0402"
Set-TextValue $ws.Range("E40") "  -18.40%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -14.52%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.857.33"
Set-TextValue $ws.Range("E42") "  -15.06%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.53"
Set-TextValue $ws.Range("E43") "  -26.74%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.256"
Set-TextValue $ws.Range("E44") "  -17.42%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0₃0570"
Set-TextValue $ws.Range("E45") "  -22.53%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -0.10%  "

# Row 47
Set-TextValue $ws.Range("B47") "InjectiveProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "25.76"
Set-TextValue $ws.Range("E47") "  -19.78%  "

# Row 48
Set-TextValue $ws.Range("B48") "ThetaToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D48") "2.31"
Set-TextValue $ws.Range("E48") "  -20.65%  "

# Row 49
Set-TextValue $ws.Range("B49") "Fetch.AI"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D49") "2.06"
Set-TextValue $ws.Range("E49") "  -19.40%  "

# Row 50
Set-TextValue $ws.Range("E50") "  -14.16%  "

# Row 51
Set-TextValue $ws.Range("B51") "Monero"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D51") "122.27"
Set-TextValue $ws.Range("E51") "  -8.42%  "

Write-Host "Applied all changes"